# CameraParam.xlsx - cam 1-8 calibration and validation
# Updates the per-camera calibration parameter values (cam3, cam4, cam5, cam6, cam8
# columns H/K/N/Q/W) for rows 6-19 (b, phi, m, k00, k10, k01, k11, k02, p00, p10,
# p01, p20, p11, p02) with freshly computed calibration/validation results.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")


# Row 6 (b)
$ws.Range("H6").Value = [double]"417.38472788313123"
$ws.Range("K6").Value = [double]"385.51517283024828"
$ws.Range("N6").Value = [double]"372.53047746730181"
$ws.Range("Q6").Value = [double]"418.61896726906627"
$ws.Range("W6").Value = [double]"414.59155703866185"

# Row 7 (phi)
$ws.Range("H7").Value = [double]"1.2318453018871041"
$ws.Range("K7").Value = [double]"1.246056901765531"
$ws.Range("N7").Value = [double]"1.2214832972117868"
$ws.Range("Q7").Value = [double]"1.223475834680041"
$ws.Range("W7").Value = [double]"1.2336373766496607"

# Row 8 (m)
$ws.Range("H8").Value = [double]"209.51888263814425"
$ws.Range("K8").Value = [double]"311.13836131048936"
$ws.Range("N8").Value = [double]"349.2570086444328"
$ws.Range("Q8").Value = [double]"213.28076740643854"
$ws.Range("W8").Value = [double]"224.3115926121603"

# Row 9 (k00)
$ws.Range("H9").Value = [double]"5.3319008124942972E-3"
$ws.Range("K9").Value = [double]"-6.1553359649351075E-2"
$ws.Range("N9").Value = [double]"7.0039324356845442E-2"
$ws.Range("Q9").Value = [double]"-0.16283449030049105"
$ws.Range("W9").Value = [double]"7.4858056619450716E-3"

# Row 10 (k10)
$ws.Range("H10").Value = [double]"1.0275586562146393"
$ws.Range("K10").Value = [double]"0.90500162007239771"
$ws.Range("N10").Value = [double]"0.93977305269963918"
$ws.Range("Q10").Value = [double]"1.0496466127945638"
$ws.Range("W10").Value = [double]"0.96426375714331347"

# Row 11 (k01)
$ws.Range("H11").Value = [double]"-1.0937444440933716E-5"
$ws.Range("K11").Value = [double]"1.5161722284629277E-4"
$ws.Range("N11").Value = [double]"-1.6066110677298942E-4"
$ws.Range("Q11").Value = [double]"3.8819431390388999E-4"
$ws.Range("W11").Value = [double]"-1.4747658351105487E-5"

# Row 12 (k11)
$ws.Range("H12").Value = [double]"-7.4551344482779513E-5"
$ws.Range("K12").Value = [double]"1.8103816404555385E-5"
$ws.Range("N12").Value = [double]"-5.2631510924545176E-5"
$ws.Range("Q12").Value = [double]"-1.4236464555665625E-4"
$ws.Range("W12").Value = [double]"-5.0456921647966996E-5"

# Row 13 (k02)
$ws.Range("H13").Value = [double]"5.0192567424011683E-9"
$ws.Range("K13").Value = [double]"-9.4655121247922585E-8"
$ws.Range("N13").Value = [double]"9.1104699471319282E-8"
$ws.Range("Q13").Value = [double]"-2.295085300724407E-7"
$ws.Range("W13").Value = [double]"8.0010843241401366E-9"

# Row 14 (p00)
$ws.Range("H14").Value = [double]"1.0103379997271895"
$ws.Range("K14").Value = [double]"1.2717505710360166"
$ws.Range("N14").Value = [double]"1.1565896041320414"
$ws.Range("Q14").Value = [double]"1.0591850835204801"
$ws.Range("W14").Value = [double]"1.1212248553798589"

# Row 15 (p10)
$ws.Range("H15").Value = [double]"-3.033460471162069E-5"
$ws.Range("K15").Value = [double]"-6.2668034949167975E-4"
$ws.Range("N15").Value = [double]"-3.8091465268874923E-4"
$ws.Range("Q15").Value = [double]"-1.5772219526196605E-4"
$ws.Range("W15").Value = [double]"-2.8565340504354723E-4"

# Row 16 (p01)
$ws.Range("H16").Value = [double]"2.7822972024428849E-3"
$ws.Range("K16").Value = [double]"0.27596444228816069"
$ws.Range("N16").Value = [double]"3.542027869425473E-2"
$ws.Range("Q16").Value = [double]"0.15397624170575072"
$ws.Range("W16").Value = [double]"-5.696386190692438E-2"

# Row 17 (p20)
$ws.Range("H17").Value = [double]"2.1201811728213964E-8"
$ws.Range("K17").Value = [double]"3.6079807671238336E-7"
$ws.Range("N17").Value = [double]"2.3157355456128215E-7"
$ws.Range("Q17").Value = [double]"1.0347318735723433E-7"
$ws.Range("W17").Value = [double]"1.6803454996444539E-7"

# Row 18 (p11)
$ws.Range("H18").Value = [double]"-7.1141887634815522E-6"
$ws.Range("K18").Value = [double]"-1.991746473730518E-4"
$ws.Range("N18").Value = [double]"-3.1215388228709629E-5"
$ws.Range("Q18").Value = [double]"-1.227085915862596E-4"
$ws.Range("W18").Value = [double]"4.1961317868103062E-5"

# Row 19 (p02)
$ws.Range("H19").Value = [double]"-0.46595272887933881"
$ws.Range("K19").Value = [double]"-0.50232028017242847"
$ws.Range("N19").Value = [double]"-0.38215182836526018"
$ws.Range("Q19").Value = [double]"-0.50731506991773123"
$ws.Range("W19").Value = [double]"-0.51722390717340161"
